$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two new helper columns (H, I) right after the "#Act" column (G) ---
# This pushes amountStud/etc. from I..R over to K..T, matching the target layout.
$ws.Columns("H:I").Insert()

# --- Insert one new row at 37 for the "totaal = 2919" summary line ---
# (pushes old rows 38/39/41 down to 39/40/42)
$ws.Rows("37:37").Insert()

# --- Column G: drop the IF(...) wrapper, just sum the three activity flags ---
$ws.Range("G2").Formula = "=B2+C2+E2"
$ws.Range("G3:G30").Formula = "=B3+C3+E3"

# --- New column H: activities beyond the first (G-1) ---
$ws.Range("H2").Formula = "=G2-1"
$ws.Range("H3:H30").Formula = "=G3-1"

# --- New column I: per-subject "student conflict" contribution (H * students) ---
$ws.Range("I2").Formula = "=H2*K2"
$ws.Range("I3:I30").Formula = "=H3*K3"

# --- New column J: total student-activity load (students * activities) ---
$ws.Range("J2").Formula = "=K2*G2"
$ws.Range("J3:J30").Formula = "=K3*G3"

# --- Row 31 totals for the new columns ---
$ws.Range("H31").Formula = "=SUM(H2:H30)"
$ws.Range("I31").Formula = "=SUM(I2:I30)"
$ws.Range("J31").Formula = "=SUM(J2:J30)"

# --- New shared strings are authored in this order so the sharedStrings.xml table
#     comes out in the same order as the target workbook (49/50/51) ---
$ws.Range("A37").Value = "totaal = 2919"
$ws.Range("A36").Value = "Slechts mogelijke score: Subjecten niet verspreid: 43*-10 = 430, 3370 studenten - 20 plaatsen in kleinste lokaaal = 3350, studentenConflict: per vak opgeteld al 1960"
$ws.Range("I32").Value = "studenten met minimaal 2 vakken"

# --- Row 32: first value of the new mini table ---
$ws.Range("J32").Value = 406

# --- Rows 33-35: remaining mini table rows (#vakken / #studenten) ---
$ws.Range("I33").Value = 3
$ws.Range("J33").Value = 241
$ws.Range("I34").Value = 4
$ws.Range("J34").Value = 116
$ws.Range("I35").Value = 5
$ws.Range("J35").Value = 38

# --- Row 36: total of the mini table ---
$ws.Range("J36").Formula = "=SUM(J31:J35)"

# --- Row 37: the new "totaal = 2919" score line ---
$ws.Range("B37").Formula = "=H31*-10-I31-(J36-20)"
$ws.Range("C37").Value = 1440
$ws.Range("D37").Formula = "=-B37+C37"

# --- Row 39: add the black-font formatting on G39:J39 (matches added cellXf/font) ---
$ws.Range("G39:J39").Font.Color = 0

# --- View tweaks: selection / top-left cell on the sheet, and window geometry on the workbook ---
$ws.Range("A36:D37").Select()
$excel.ActiveWindow.ScrollRow = 4

$wb.Windows.Item(1).Left = 60
$wb.Windows.Item(1).Width = 25480
